$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timetable cell contents (B2:F13) to reflect the newly generated schedule
$ws.Range("B2").Value = '{0: sala nr 1 | Dominik Kaczor | Informatyka}'
$ws.Range("C2").Value = '{}'
$ws.Range("D2").Value = '{}'
$ws.Range("E2").Value = '{0: sala nr 9 | Piotr Wójcik | Biologia}'
$ws.Range("F2").Value = '{}'

$ws.Range("B3").Value = '{0: sala nr 8 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("C3").Value = '{}'
$ws.Range("D3").Value = '{}'
$ws.Range("E3").Value = '{0: sala nr 7 | Zofia Wiśniewska | Wychowanie fizyczne}'
$ws.Range("F3").Value = '{0: sala nr 2 | Jan Nowak | Język polski}'

$ws.Range("B4").Value = '{0: sala nr 6 | Paweł Lewandowski | Matematyka}'
$ws.Range("C4").Value = '{}'
$ws.Range("D4").Value = '{}'
$ws.Range("E4").Value = '{0: sala nr 4 | Piotr Wójcik | Biologia}'
$ws.Range("F4").Value = '{0: sala nr 1 | Mateusz Kowalski | Język niemiecki}'

$ws.Range("B5").Value = '{0: sala nr 11 | Lena Kowalska | Język angielski}'
$ws.Range("C5").Value = '{}'
$ws.Range("D5").Value = '{}'
$ws.Range("E5").Value = '{}'
$ws.Range("F5").Value = '{0: sala nr 9 | Dominik Kaczor | Informatyka}'

$ws.Range("B6").Value = '{}'
$ws.Range("C6").Value = '{0: sala nr 10 | Paweł Lewandowski | Matematyka}'
$ws.Range("D6").Value = '{}'
$ws.Range("E6").Value = '{}'
$ws.Range("F6").Value = '{0: sala nr 5 | Natalia Szymańska | Geografia}'

$ws.Range("B7").Value = '{}'
$ws.Range("C7").Value = '{0: sala nr 6 | Jan Nowak | Język polski}'
$ws.Range("D7").Value = '{}'
$ws.Range("E7").Value = '{0: sala nr 5 | Karolina Kamińska | Chemia}'
$ws.Range("F7").Value = '{0: sala nr 7 | Paweł Lewandowski | Matematyka}'

$ws.Range("B8").Value = '{}'
$ws.Range("C8").Value = '{0: sala nr 1 | Lena Kowalska | Język angielski}'
$ws.Range("D8").Value = '{}'
$ws.Range("E8").Value = '{0: sala nr 7 | Natalia Szymańska | Geografia}'
$ws.Range("F8").Value = '{}'

$ws.Range("B9").Value = '{}'
$ws.Range("C9").Value = '{0: sala nr 6 | Katarzyna Mazur | Fizyka}'
$ws.Range("D9").Value = '{0: sala nr 4 | Katarzyna Mazur | Fizyka}'
$ws.Range("E9").Value = '{0: sala nr 10 | Mateusz Kowalski | Język niemiecki}'
$ws.Range("F9").Value = '{}'

$ws.Range("B10").Value = '{}'
$ws.Range("C10").Value = '{}'
$ws.Range("D10").Value = '{0: sala nr 9 | Dominik Kaczor | Informatyka}'
$ws.Range("E10").Value = '{0: sala nr 11 | Paweł Lewandowski | Matematyka}'
$ws.Range("F10").Value = '{0: sala nr 5 | Paweł Lewandowski | Matematyka}'

$ws.Range("B11").Value = '{}'
$ws.Range("C11").Value = '{}'
$ws.Range("D11").Value = '{0: sala nr 5 | Jan Nowak | Język polski}'
$ws.Range("E11").Value = '{}'
$ws.Range("F11").Value = '{0: sala nr 8 | Zofia Wiśniewska | Wychowanie fizyczne}'

$ws.Range("B12").Value = '{}'
$ws.Range("C12").Value = '{}'
$ws.Range("D12").Value = '{0: sala nr 4 | Katarzyna Mazur | Fizyka}'
$ws.Range("E12").Value = '{0: sala nr 9 | Lena Kowalska | Język angielski}'
$ws.Range("F12").Value = '{0: sala nr 10 | Dominik Kaczor | Informatyka}'

$ws.Range("B13").Value = '{}'
$ws.Range("C13").Value = '{}'
$ws.Range("D13").Value = '{0: sala nr 9 | Karolina Kamińska | Chemia}'
$ws.Range("E13").Value = '{0: sala nr 5 | Dominik Kaczor | Informatyka}'
$ws.Range("F13").Value = '{}'

# Update column widths (B, C, D, E) to match the new layout
$ws.Columns.Item(2).ColumnWidth = 54.83
$ws.Columns.Item(3).ColumnWidth = 47.83
$ws.Columns.Item(4).ColumnWidth = 44.83
$ws.Columns.Item(5).ColumnWidth = 54.83
